# Feria Lagunitas de Puerto Montt - Haba
# Insert a new weekly data row (row 64) in the sheet, pushing all
# subsequent records down by one row (old row 64 -> new row 65, ...,
# old row 99 -> new row 100), and populate the new row 64 with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64:99 down to 65:100, inserting a blank row at 64.
$ws.Rows(64).Insert()

# Fill the new row 64 with the new weekly record.
$ws.Cells.Item(64, 1).Value  = 4
$ws.Cells.Item(64, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value  = "Los Lagos"
$ws.Cells.Item(64, 4).Value  = 44813
$ws.Cells.Item(64, 5).Value  = 10
$ws.Cells.Item(64, 6).Value  = 100112026
$ws.Cells.Item(64, 7).Value  = "Haba"
$ws.Cells.Item(64, 8).Value  = "Sin especificar"
$ws.Cells.Item(64, 9).Value  = "Primera"
$ws.Cells.Item(64, 10).Value = 80
$ws.Cells.Item(64, 11).Value = 15000
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 13).Value = 15000
$ws.Cells.Item(64, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(64, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 16).Value = 600
$ws.Cells.Item(64, 17).Value = 25
$ws.Cells.Item(64, 18).Value = "Hortaliza"
